# Actualización automática hashcode vie ago 30 02:20:30 CEST 2019
#
# Updates 14 existing hashcode values in column B, then appends 12 new
# rows (id in column A, "nuevo" literal in column B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update existing hash values (row -> new hash) ---------------------
$updates = @{
    44  = "a2cfcbfef9b7b4aed5ed06cdf76e820f"
    74  = "9555bf74da8a390313ded720eb47dce7"
    89  = "160ee88f449d69ffbf488ebe9d2dcc44"
    99  = "ec5bd2a050b8a245967e920be6cdaaa2"
    110 = "4050bd447a74401c61ea746f9711d4fc"
    123 = "451e5ab82ec5458c7ce53697b094da82"
    161 = "9bb4c7968671c6ffbee5b3db18131f17"
    168 = "36c8cd53ba8a46717318adc0a51706b1"
    278 = "4f4e6e1d7f91885a3a4f184b8ac396e3"
    345 = "183913fecc02620ae6913e0667b17656"
    768 = "8a866f38cea4d509d812189b47eef642"
    816 = "1951623ae9020a139ec3467817acc2ab"
    825 = "76fb08e3968f1341beee8c4d704ab1a6"
    827 = "fe391b223dd9b3e7fc6a5f6ebd9890a3"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}

# --- 2) Append 12 new rows at the bottom -----------------------------------
$newIds = @(
    "901-010202-9010102021TM",
    "901-010202-9010102021TP",
    "901-010202-9010102023TC",
    "901-010202-9010102022A",
    "901-010202-9010102021A",
    "901-010202-9010102022TC",
    "901-010202-9010102023TM",
    "901-010202-9010102023TP",
    "901-010202-9010102023A",
    "901-010202-9010102021TC",
    "901-010202-9010102022TM",
    "901-010202-9010102022TP"
)

$startRow = 963
for ($i = 0; $i -lt $newIds.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newIds[$i]
    $ws.Cells.Item($r, 2).Value = "nuevo"
}
